# DPP master file export: the 2nd table's "SESSION 17" column header was
# wrapped onto two lines ("SESSION 17" / "10/01/2019"); flatten it to a
# single line to match the rest of the session headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined")

$ws.Range("X1").Value = "SESSION 17 10/01/2019"

$ws.Range("X1").Select()
